# Fill in column Q ("uds. Objetivo semana pasada") with the values
# reported for last week's target units, per the commit:
# "Hemos conseguido rellenar la columna uds. objetivo semana pasada"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_8")

$updates = @{
    3  = 4
    5  = 16
    8  = 2
    9  = 2
    10 = 2
    12 = 2
    13 = 1
    14 = 1
    16 = 2
    19 = 3
    20 = 1
    21 = 1
    24 = 2
    25 = 3
    26 = 6
    31 = 2
    33 = 4
    34 = 2
    35 = 2
    36 = 1
    40 = 5
    41 = 22
    42 = 40
    43 = 4
    44 = 14
    45 = 6
    46 = 4
    50 = 2
    51 = 3
    52 = 4
    56 = 12
    57 = 1
    58 = 7
    59 = 1
    60 = 9
    63 = 1
    64 = 7
    65 = 7
    66 = 141
    67 = 5
    68 = 3
    70 = 2
    71 = 3
    72 = 1
    73 = 2
    80 = 1
    81 = 2
    84 = 2
    85 = 5
    86 = 1
    87 = 2
    88 = 1
    90 = 1
    92 = 3
    93 = 8
}

foreach ($row in $updates.Keys) {
    $ws.Range("Q$row").Value = $updates[$row]
}
